$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1: copy style from H1 (reuses existing style index 1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$excel.CutCopyMode = 0

$iVals = @(8,7,9,9,9,8,9,9,9,7,8,9,9,9,9,8,9,9,9,9,7,7,7,7,8,9,8,9,9,7,8,9,8,7,9,9,8,7,8,8,7,9,8,8,7,8,7,8,8,8,8,8,8,6,7,7,6,9,9,8,8,6,8,6,5,3,6,6,5)
$jVals = @(8,7,10,9,9,8,9,9,9,7,8,9,9,9,9,9,9,9,9,9,7,7,7,8,9,9,9,9,9,8,8,9,8,7,9,9,8,7,8,8,8,9,8,8,7,8,7,8,8,8,8,8,8,7,7,7,6,9,9,8,8,6,8,6,5,3,6,6,5)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$k]
    $ws.Cells.Item($r, 10).Value = $jVals[$k]
}
